$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.650.29"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "3.580.45"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.72"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.40"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.491"
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").Value = "4.188.31"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.01"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "3.604.10"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "66.655.02"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.116"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.44"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.01"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "431.97"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.619"
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.19"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").Value = "3.722.68"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.06"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.19"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "3.577.35"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.42"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("E33").Value = "  -3.47%  "
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.84"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.62"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "173.60"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.93"
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  +5.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.19"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.15"
$ws.Range("E47").Value = "  -3.43%  "
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.63"
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.941"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("E51").Value = "  -1.42%  "
